$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($addr in @("D5", "D6", "D7", "D9", "D10", "D11", "D12", "D13", "D16", "D17", "D20", "D21", "D23", "D24", "D25", "D26", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D36", "D37", "D39", "D40", "D41", "D42", "D45", "D46", "D47", "D48", "D49", "D50", "D51")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "42.302.47"
$ws.Range("E2").Value = "  -7.16%  "
$ws.Range("D3").Value = "2.198.29"
$ws.Range("E3").Value = "  -7.67%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "307.72"
$ws.Range("E5").Value = "  -2.27%  "
$ws.Range("D6").Value = "96.64"
$ws.Range("E6").Value = "  -13.41%  "
$ws.Range("D7").Value = "0.567"
$ws.Range("E7").Value = "  -9.96%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "0.547"
$ws.Range("E9").Value = "  -11.41%  "
$ws.Range("D10").Value = "36.04"
$ws.Range("E10").Value = "  -12.42%  "
$ws.Range("D11").Value = "53.45"
$ws.Range("E11").Value = "  -4.76%  "
$ws.Range("D12").Value = "0.0817"
$ws.Range("E12").Value = "  -11.49%  "
$ws.Range("D13").Value = "7.45"
$ws.Range("E13").Value = "  -13.00%  "
$ws.Range("E14").Value = "  -4.92%  "
$ws.Range("D15").Value = "2.538.08"
$ws.Range("E15").Value = "  -7.66%  "
$ws.Range("D16").Value = "0.847"
$ws.Range("E16").Value = "  -14.22%  "
$ws.Range("D17").Value = "13.73"
$ws.Range("E17").Value = "  -12.02%  "
$ws.Range("D18").Value = "2.221.97"
$ws.Range("E18").Value = "  -7.57%  "
$ws.Range("D19").Value = "42.253.65"
$ws.Range("E19").Value = "  -7.16%  "
$ws.Range("D20").Value = "13.84"
$ws.Range("E20").Value = "  +5.52%  "
$ws.Range("D21").Value = "6.44"
$ws.Range("E21").Value = "  -12.34%  "
$ws.Range("D22").Value = "0.0₃0931"
$ws.Range("E22").Value = "  -12.97%  "
$ws.Range("B23").Value = "PancakeSwap"
$ws.Range("C23").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D23").Value = "3.15"
$ws.Range("E23").Value = "  -9.61%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "64.28"
$ws.Range("E24").Value = "  -13.01%  "
$ws.Range("D25").Value = "230.84"
$ws.Range("E25").Value = "  -12.01%  "
$ws.Range("D26").Value = "2.07"
$ws.Range("E26").Value = "  -10.05%  "
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("D28").Value = "9.96"
$ws.Range("E28").Value = "  -10.66%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "2.13"
$ws.Range("E29").Value = "  -9.58%  "
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").Value = "6.35"
$ws.Range("E30").Value = "  -15.17%  "
$ws.Range("D31").Value = "20.18"
$ws.Range("E31").Value = "  -10.54%  "
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").Value = "156.73"
$ws.Range("E32").Value = "  -8.17%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "0.0855"
$ws.Range("E33").Value = "  -11.47%  "
$ws.Range("D34").Value = "32.93"
$ws.Range("E34").Value = "  -13.80%  "
$ws.Range("E35").Value = "  -9.60%  "
$ws.Range("D36").Value = "3.17"
$ws.Range("E36").Value = "  +5.83%  "
$ws.Range("D37").Value = "0.120"
$ws.Range("E37").Value = "  -8.22%  "
$ws.Range("E38").Value = "  -9.08%  "
$ws.Range("D39").Value = "1.79"
$ws.Range("E39").Value = "  +3.33%  "
$ws.Range("D40").Value = "0.102"
$ws.Range("E40").Value = "  -12.99%  "
$ws.Range("D41").Value = "3.45"
$ws.Range("E41").Value = "  -13.76%  "
$ws.Range("D42").Value = "0.0312"
$ws.Range("E42").Value = "  -12.54%  "
$ws.Range("E43").Value = "  +0.13%  "
$ws.Range("D44").Value = "1.756.92"
$ws.Range("E44").Value = "  +6.46%  "
$ws.Range("D45").Value = "86.85"
$ws.Range("E45").Value = "  -15.29%  "
$ws.Range("D46").Value = "11.74"
$ws.Range("E46").Value = "  -11.59%  "
$ws.Range("D47").Value = "0.201"
$ws.Range("E47").Value = "  -14.95%  "
$ws.Range("D48").Value = "74.94"
$ws.Range("E48").Value = "  -11.10%  "
$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D49").Value = "5.21"
$ws.Range("E49").Value = "  -6.45%  "
$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").Value = "59.11"
$ws.Range("E50").Value = "  -15.96%  "
$ws.Range("D51").Value = "8.38"
$ws.Range("E51").Value = "  -10.77%  "
